$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '93.884.38'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.062.36'
$ws.Range("E3").Value = '  -2.21%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.19'
$ws.Range("E5").Value = '  -4.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '604.56'
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.09'
$ws.Range("E7").Value = '  -2.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.373'
$ws.Range("E8").Value = '  -8.44%  '
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  +6.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.054.86'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("E12").Value = '  -4.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '93.441.56'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000238'
$ws.Range("E14").Value = '  -6.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '33.33'
$ws.Range("E15").Value = '  -4.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.25'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.620.70'
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.051.35'
$ws.Range("E18").Value = '  -3.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.48'
$ws.Range("E19").Value = '  -8.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.23'
$ws.Range("E20").Value = '  -4.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.63'
$ws.Range("E21").Value = '  -3.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '435.59'
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.69'
$ws.Range("E23").Value = '  -7.87%  '
$ws.Range("E24").Value = '  -10.66%  '
$ws.Range("E25").Value = '  +5.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '5.44'
$ws.Range("E26").Value = '  -7.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '84.23'
$ws.Range("E27").Value = '  -4.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.57'
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.196.69'
$ws.Range("E29").Value = '  -3.10%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  +5.27%  '
$ws.Range("E32").Value = '  +12.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.175'
$ws.Range("E33").Value = '  +2.57%  '
$ws.Range("E34").Value = '  -11.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.94'
$ws.Range("E35").Value = '  -4.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.58'
$ws.Range("E36").Value = '  -7.01%  '
$ws.Range("E37").Value = '  -4.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.10'
$ws.Range("E38").Value = '  -5.07%  '
$ws.Range("E39").Value = '  -2.10%  '
$ws.Range("E40").Value = '  +3.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.74'
$ws.Range("E41").Value = '  -3.68%  '
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '458.26'
$ws.Range("E43").Value = '  -7.59%  '
$ws.Range("E44").Value = '  -6.30%  '
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.08'
$ws.Range("E46").Value = '  -11.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '159.91'
$ws.Range("E47").Value = '  -2.43%  '
$ws.Range("E48").Value = '  -7.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.658'
$ws.Range("E49").Value = '  -5.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.62'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("E51").Value = '  -0.01%  '
